$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FTV")

# Row 6: "Change in inventories" - update B6:G6
$ws.Range("B6").Value = 5500000.0
$ws.Range("C6").Value = -7300000.0
$ws.Range("D6").Value = 629800000.0
$ws.Range("E6").Value = 634700000.0
$ws.Range("F6").Value = 726100000.0
$ws.Range("G6").Value = 709400000.0

# Row 8: "Change in payables and accrued liability" - update B8:G8
$ws.Range("B8").Value = 1910400000.0
$ws.Range("C8").Value = 2736000000.0
$ws.Range("D8").Value = 3021000000.0
$ws.Range("E8").Value = 2375100000.0
$ws.Range("F8").Value = 1651700000.0
$ws.Range("G8").Value = 893500000.0
